$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.731.08"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").Value = "3.106.27"
$ws.Range("E3").Value = "  +3.64%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "559.49"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "143.44"
$ws.Range("E6").Value = "  +9.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.100.36"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").Value = "7.15"
$ws.Range("E10").Value = "  +19.51%  "
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("E12").Value = "  +4.45%  "
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "35.41"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "3.607.00"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "64.710.59"
$ws.Range("E16").Value = "  +4.88%  "
$ws.Range("D17").Value = "3.107.24"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "484.07"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "13.83"
$ws.Range("E21").Value = "  +5.32%  "
$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  +9.89%  "
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  +10.98%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  +9.18%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "26.15"
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("D33").Value = "2.46"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").Value = "55.34"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").Value = "464.66"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("E38").Value = "  +7.60%  "
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("D40").Value = "3.018.39"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("D42").Value = "8.29"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +15.86%  "
$ws.Range("D44").Value = "28.34"
$ws.Range("E44").Value = "  +11.87%  "
$ws.Range("D45").Value = "0.262"
$ws.Range("E45").Value = "  +8.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  +8.78%  "
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("D49").Value = "118.92"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("E50").Value = "  +7.14%  "
$ws.Range("E51").Value = "  +3.27%  "
